$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sensor label
$ws.Range("A1").Value = "Sensor 3"

# Update the range labels (B column)
$ws.Range("B3").Value = "Dato Min: 1"
$ws.Range("B4").Value = "Dato Actual: 2000"
$ws.Range("B5").Value = "Dato Max: 15000"

# Update the sensor history data values (A column)
$ws.Range("A3").Value = 10
$ws.Range("A4").Value = 1600
$ws.Range("A5").Value = 3300
$ws.Range("A6").Value = 5000
$ws.Range("A7").Value = 1000
$ws.Range("A8").Value = 20000
$ws.Range("A9").Value = 2000

# Remove row 10 entirely (was A10 = 120)
$ws.Rows("10").Delete()
